# edit.ps1
# Applies stock-movement data corrections to the "CryCompanywiseStockReport_1"
# worksheet: several rows had their Batch No. (B), Sale Rate (E), Qty (F) and
# Amount (G) values out of order relative to their neighbours. This script
# writes the corrected values directly into the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B146" = 57756
    "F146" = -100
    "G146" = -6644
    "B147" = 53925
    "E147" = 79.37
    "F147" = 1
    "G147" = 66.44
    "B148" = 64350
    "E148" = 70.63
    "F148" = 2
    "G148" = 132.88
    "B163" = 64329
    "E163" = 128.32
    "F163" = 4
    "G163" = 482.76
    "B164" = 57552
    "E164" = 136.86
    "F164" = -5
    "G164" = -603.45
    "B246" = 48706
    "E246" = 39.8
    "F246" = -144
    "G246" = -4795.2
    "B247" = 64973
    "E247" = 35.4
    "F247" = 128
    "G247" = 4262.4
    "B294" = 63531
    "F294" = 80
    "G294" = 11478.4
    "B295" = 63571
    "F295" = 10
    "G295" = 1434.8
    "B299" = 63510
    "E299" = 50.66
    "F299" = 152
    "G299" = 7241.28
    "B300" = 55356
    "E300" = 54.04
    "F300" = -158
    "G300" = -7527.12
    "B315" = 63560
    "E315" = 134.87
    "F315" = 1
    "G315" = 126.86
    "B316" = 60325
    "E316" = 151.57
    "F316" = -102
    "G316" = -12939.72
    "B356" = 63681
    "E356" = 23.84
    "F356" = 0
    "G356" = 0
    "B357" = 31930
    "E357" = 26.8
    "F357" = -62
    "G357" = -1390.04
    "B472" = 45695
    "E472" = 23.58
    "F472" = -36
    "G472" = -710.28
    "B473" = 64915
    "E473" = 20.98
    "F473" = 0
    "G473" = 0
    "B479" = 64927
    "E479" = 17.26
    "F479" = 253
    "G479" = 4103.66
    "B480" = 45718
    "E480" = 19.38
    "F480" = -294
    "G480" = -4768.68
    "B564" = 64810
    "E564" = 291.22
    "F564" = 6
    "G564" = 1643.52
    "B565" = 53319
    "E565" = 310.64
    "F565" = -6
    "G565" = -1643.52
    "B596" = 60022
    "E596" = 37.22
    "F596" = -113
    "G596" = -3709.79
    "B597" = 64830
    "E597" = 34.9
    "F597" = 113
    "G597" = 3709.79
    "B705" = 63150
    "D705" = 75.68000000000001
    "E705" = 80.45
    "F705" = 112
    "G705" = 8476.16
    "B706" = 61428
    "D706" = 69.16
    "E706" = 73.52
    "F706" = 1
    "G706" = 69.16
    "B732" = 65079
    "F732" = 21
    "G732" = 858.27
    "B733" = 65362
    "F733" = 75
    "G733" = 3065.25
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
